$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vuelta 1")

# Remove the duplicate "MULTITRANS RR" row (old row 116); this shifts all
# subsequent rows up by one and shrinks the used range from N121 to N120.
$ws.Rows(116).Delete()

# Row 115 becomes a new order: APM-857 / UL&PH INVERSIONES SAC
$ws.Range("B115").Value = "APM-857"
$ws.Range("C115").Value = "UL&PH INVERSIONES SAC"
$ws.Range("D115").Value = 7600
$ws.Range("E115").Value = 8
$ws.Range("N115").Value = 1000
